$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Header row stays the same except C1 changes from "cond" to "condition"
$ws.Range("C1").Value = "condition"

# Row 2: coles / specialpricecat product / specialpricecat==012 & product==SJ300
$ws.Range("A2").Value = "coles"
$ws.Range("B2").Value = "specialpricecat product"
$ws.Range("C2").Value = "specialpricecat==012 & product==SJ300"

# Row 3: woolworths / specialpricecat product / specialpricecat==010 & product==SJ300
$ws.Range("A3").Value = "woolworths"
$ws.Range("B3").Value = "specialpricecat product"
$ws.Range("C3").Value = "specialpricecat==010 & product==SJ300"

# Row 4: shop / specialpricecat product / specialpricecat==092 & product==SJ300
$ws.Range("A4").Value = "shop"
$ws.Range("B4").Value = "specialpricecat product"
$ws.Range("C4").Value = "specialpricecat==092 & product==SJ300"

# Row 5: indies / cat product / cat==88 & product==SJ300
$ws.Range("A5").Value = "indies"
$ws.Range("B5").Value = "cat product"
$ws.Range("C5").Value = "cat==88 & product==SJ300"

# Row 6: distributors / cat code product / cat==81 & product==SJ300
$ws.Range("A6").Value = "distributors"
$ws.Range("B6").Value = "cat code product"
$ws.Range("C6").Value = "cat==81 & product==SJ300"

# Update selection to C1 as per diff
$ws.Range("C1").Select()
